$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137, shifting existing rows 137:168 down to 138:169
$ws.Rows("137:137").Insert()

# Populate the newly inserted row 137 with the new weekly record
$ws.Range("A137").Value = 4
$ws.Range("B137").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C137").Value = "Los Lagos"
$ws.Range("D137").Value = 44508
$ws.Range("E137").Value = 10
$ws.Range("F137").Value = "Fruta"
$ws.Range("G137").Value = 100102
$ws.Range("H137").Value = "Cítricos"
$ws.Range("I137").Value = 100102006
$ws.Range("J137").Value = "Pomelo"
$ws.Range("K137").Value = "Start Ruby"
$ws.Range("L137").Value = "Primera"
$ws.Range("M137").Value = 60
$ws.Range("N137").Value = 11000
$ws.Range("O137").Value = 12000
$ws.Range("P137").Value = 11500
$ws.Range("Q137").Value = "$/caja 14 kilos empedrada"
$ws.Range("R137").Value = "Región de O'Higgins"
$ws.Range("S137").Value = 821
$ws.Range("T137").Value = 14
